# Update the hourly forecast report figures ("align chart starting from 8
# for real-time line"). Only the numeric measure columns (C,D,E,G,H,I,J and
# one B) change; row/category labels in column A stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

Set-Row 2  @{ C = 1.237; D = 1.242; G = 22; I = 1315; J = -99.90555133079849 }
Set-Row 3  @{ C = 407;   D = 411;          I = 380;  J = 8.157894736842097 }
Set-Row 4  @{ C = 19;    D = 19;           I = 18;   J = 5.555555555555558 }
Set-Row 5  @{ C = 265;   D = 275;          I = 323;  J = -14.86068111455109 }
Set-Row 6  @{ C = 82;    D = 83;           I = 89;   J = -6.741573033707871 }
Set-Row 7  @{ C = 83;    D = 85;           I = 127;  J = -33.07086614173228 }
Set-Row 8  @{ C = 251;   D = 258;  G = 7;  I = 451;  J = -42.79379157427938 }
Set-Row 9  @{ C = 39;    D = 38;           I = 59;   J = -35.59322033898306 }
Set-Row 10 @{ C = 245;   D = 279;  E = 37; I = 909;  J = -69.3069306930693 }
Set-Row 11 @{ C = 308;   D = 310;  G = 7;  I = 407;  J = -23.83292383292384 }
Set-Row 12 @{ B = 1; C = 439; D = 633; H = 118; I = 1117; J = -43.3303491495076 }
Set-Row 13 @{ C = 30;    D = 30;           I = 633;  J = -95.260663507109 }
Set-Row 14 @{ C = 596;   D = 716;  G = 11; H = 74; I = 1027; J = -30.2823758519961 }
Set-Row 15 @{ C = 237;   D = 267;          I = 306;  J = -12.74509803921569 }
Set-Row 16 @{ C = 148;   D = 156;  G = 5;  I = 164;  J = -4.878048780487809 }
Set-Row 17 @{ C = 5;     D = 8;            J = 100 }
Set-Row 18 @{ C = 21;    D = 21;           I = 22;   J = -4.545454545454541 }
Set-Row 19 @{ C = 14;    D = 15;  G = 1;   I = 89;   J = -83.14606741573034 }
